$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of trade data (rows 6-8)
$ws.Range("A6").Value = 9990.1
$ws.Range("B6").Value = 9959.23
$ws.Range("C6").Value = 107.89
$ws.Range("D6").Value = 108.22
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = 0.31
$ws.Range("G6").Value = 42613.766574074078
$ws.Range("H6").Value = $true

$ws.Range("A7").Value = 9994.1
$ws.Range("B7").Value = 9990.1
$ws.Range("C7").Value = 107.17
$ws.Range("D7").Value = 107.21
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = 0.04
$ws.Range("G7").Value = 42614.673807870371
$ws.Range("H7").Value = $true

$ws.Range("A8").Value = 9992.1
$ws.Range("B8").Value = 9994.1
$ws.Range("C8").Value = 107.04
$ws.Range("D8").Value = 107.02
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = -0.02
$ws.Range("G8").Value = 42615.752905092595
$ws.Range("H8").Value = $false
